$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename CDE mole-concentration observations (rows unaffected structurally) ---
$ws.Range("A14").Value = "cde-chloride-in-serum-or-plasma-mole-concentration"
$ws.Range("B14").Value = "CDE Chloride In Serum Or Plasma Mole Concentration"
$ws.Range("A15").Value = "cde-creatinine-in-serum-or-plasma-mole-concentration"
$ws.Range("B15").Value = "CDE Creatinine In Serum Or Plasma Mole Concentration"
$ws.Range("A27").Value = "cde-glucose-in-serum-or-plasma-mole-concentration"
$ws.Range("B27").Value = "CDE Glucose In Serum Or Plasma Mole Concentration"
$ws.Range("A29").Value = "cde-hemoglobin-in-blood-mole-concentration"
$ws.Range("B29").Value = "CDE Hemoglobin In Blood Mole Concentration"
$ws.Range("A32").Value = "cde-magnesium-in-serum-or-plasma-mole-concentration"
$ws.Range("B32").Value = "CDE Magnesium In Serum Or Plasma Mole Concentration"
$ws.Range("A33").Value = "cde-phosphate-in-serum-or-plasma-mole-concentration"
$ws.Range("B33").Value = "CDE Phosphate In Serum Or Plasma Mole Concentration"
$ws.Range("A34").Value = "cde-potassium-in-serum-or-plasma-mole-concentration"
$ws.Range("B34").Value = "CDE Potassium In Serum Or Plasma Mole Concentration"
$ws.Range("A39").Value = "cde-sodium-in-serum-or-plasma-mole-concentration"
$ws.Range("B39").Value = "CDE Sodium In Serum Or Plasma Mole Concentration"

# --- Rewrite classification rows 49-70 with the new CLS OBO taxonomy ---
# Row 49
$ws.Range("A49").Value = "cls-obo-cmo-blood-calcium-level"
$ws.Range("B49").Value = "CLS CMO Blood Calcium Level"
$ws.Range("C49").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000063, null#CMO_0000502"
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = ""
$ws.Range("F49").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G49").Value = "dateTime, Period, Timing, instant"
$ws.Range("H49").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I49").Value = "optional"
$ws.Range("J49").Value = ""
$ws.Range("K49").Value = ""

# Row 50
$ws.Range("A50").Value = "cls-obo-cmo-blood-chloride-level"
$ws.Range("B50").Value = "CLS OBO CMO Blood Chloride Level"
$ws.Range("C50").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000063, null#CMO_0000497"
$ws.Range("D50").Value = ""
$ws.Range("E50").Value = ""
$ws.Range("F50").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G50").Value = "dateTime, Period, Timing, instant"
$ws.Range("H50").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I50").Value = "optional"
$ws.Range("J50").Value = ""
$ws.Range("K50").Value = ""

# Row 51
$ws.Range("A51").Value = "cls-obo-cmo-blood-glucose-level"
$ws.Range("B51").Value = "CLS OBO CMO Blood Glucose Level"
$ws.Range("C51").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000046"
$ws.Range("D51").Value = ""
$ws.Range("E51").Value = ""
$ws.Range("F51").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G51").Value = "dateTime, Period, Timing, instant"
$ws.Range("H51").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I51").Value = "optional"
$ws.Range("J51").Value = ""
$ws.Range("K51").Value = ""

# Row 52
$ws.Range("A52").Value = "cls-obo-cmo-blood-magnesium-level"
$ws.Range("B52").Value = "CLS OBO CMO Blood Magnesium Level"
$ws.Range("C52").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000063, null#CMO_0000505"
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = ""
$ws.Range("F52").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G52").Value = "dateTime, Period, Timing, instant"
$ws.Range("H52").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I52").Value = "optional"
$ws.Range("J52").Value = ""
$ws.Range("K52").Value = ""

# Row 53
$ws.Range("A53").Value = "cls-obo-cmo-blood-phosphate-level"
$ws.Range("B53").Value = "CLS OBO CMO Blood Phosphate Level"
$ws.Range("C53").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000063, null#CMO_0000504"
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = ""
$ws.Range("F53").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G53").Value = "dateTime, Period, Timing, instant"
$ws.Range("H53").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I53").Value = "optional"
$ws.Range("J53").Value = ""
$ws.Range("K53").Value = ""

# Row 54
$ws.Range("A54").Value = "cls-obo-cmo-blood-potassium-level"
$ws.Range("B54").Value = "CLS OBO CMO Blood Potassium Level"
$ws.Range("C54").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000063, null#CMO_0000496"
$ws.Range("D54").Value = ""
$ws.Range("E54").Value = ""
$ws.Range("F54").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G54").Value = "dateTime, Period, Timing, instant"
$ws.Range("H54").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I54").Value = "optional"
$ws.Range("J54").Value = ""
$ws.Range("K54").Value = ""

# Row 55
$ws.Range("A55").Value = "cls-obo-cmo-blood-sodium-level"
$ws.Range("B55").Value = "CLS OBO CMO Blood Sodium Level"
$ws.Range("C55").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000063, null#CMO_0000499"
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = ""
$ws.Range("F55").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G55").Value = "dateTime, Period, Timing, instant"
$ws.Range("H55").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I55").Value = "optional"
$ws.Range("J55").Value = ""
$ws.Range("K55").Value = ""

# Row 56
$ws.Range("A56").Value = "cls-obo-cmo-hematocrit-measurement"
$ws.Range("B56").Value = "CLS OBO CMO Hematocrit Measurement"
$ws.Range("C56").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000036, null#CMO_0000508"
$ws.Range("D56").Value = ""
$ws.Range("E56").Value = ""
$ws.Range("F56").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G56").Value = "dateTime, Period, Timing, instant"
$ws.Range("H56").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I56").Value = "optional"
$ws.Range("J56").Value = ""
$ws.Range("K56").Value = ""

# Row 57
$ws.Range("A57").Value = "cls-obo-cmo-hemoglobin-measurement"
$ws.Range("B57").Value = "CLS OBO CMO Hemoglobin Measurement"
$ws.Range("C57").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000036, null#CMO_0000508"
$ws.Range("D57").Value = ""
$ws.Range("E57").Value = ""
$ws.Range("F57").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G57").Value = "dateTime, Period, Timing, instant"
$ws.Range("H57").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I57").Value = "optional"
$ws.Range("J57").Value = ""
$ws.Range("K57").Value = ""

# Row 58
$ws.Range("A58").Value = "cls-obo-cmo-platelet-measurement"
$ws.Range("B58").Value = "CLS OBO CMO Platelet Measurement"
$ws.Range("C58").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000036, null#CMO_0000921"
$ws.Range("D58").Value = ""
$ws.Range("E58").Value = ""
$ws.Range("F58").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G58").Value = "dateTime, Period, Timing, instant"
$ws.Range("H58").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I58").Value = "optional"
$ws.Range("J58").Value = ""
$ws.Range("K58").Value = ""

# Row 59
$ws.Range("A59").Value = "cls-obo-cmo-red-blood-cell-measurement"
$ws.Range("B59").Value = "CLS OBO CMO Red Blood Cell Measurement"
$ws.Range("C59").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000036, null#CMO_0001356"
$ws.Range("D59").Value = ""
$ws.Range("E59").Value = ""
$ws.Range("F59").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G59").Value = "dateTime, Period, Timing, instant"
$ws.Range("H59").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I59").Value = "optional"
$ws.Range("J59").Value = ""
$ws.Range("K59").Value = ""

# Row 60
$ws.Range("A60").Value = "cls-obo-cmo-white-blood-cell-measurement"
$ws.Range("B60").Value = "CLS OBO CMO White Blood Cell Measurement"
$ws.Range("C60").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000036, null#CMO_0002341"
$ws.Range("D60").Value = ""
$ws.Range("E60").Value = ""
$ws.Range("F60").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G60").Value = "dateTime, Period, Timing, instant"
$ws.Range("H60").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I60").Value = "optional"
$ws.Range("J60").Value = ""
$ws.Range("K60").Value = ""

# Row 61
$ws.Range("A61").Value = "cls-obo-ncit-date-of-diagnosis"
$ws.Range("B61").Value = "CLS OBO NCIT Date Of Diagnosis"
$ws.Range("C61").Value = "Observation Category Codes#social-history, null#NCIT_C20189, null#NCIT_C41009, null#NCIT_C21514, null#NCIT_C25164, null#NCIT_C164339"
$ws.Range("D61").Value = ""
$ws.Range("E61").Value = ""
$ws.Range("F61").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G61").Value = "dateTime, Period, Timing, instant"
$ws.Range("H61").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I61").Value = "optional"
$ws.Range("J61").Value = ""
$ws.Range("K61").Value = ""

# Row 62
$ws.Range("A62").Value = "cls-obo-ncit-visit-date"
$ws.Range("B62").Value = "CLS OBO NCIT Visit Date"
$ws.Range("C62").Value = "Observation Category Codes#laboratory, null#NCIT_C20189, null#NCIT_C41009, null#NCIT_C21514, null#NCIT_C25164, null#NCIT_C83031"
$ws.Range("D62").Value = ""
$ws.Range("E62").Value = ""
$ws.Range("F62").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G62").Value = "dateTime, Period, Timing, instant"
$ws.Range("H62").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I62").Value = "optional"
$ws.Range("J62").Value = ""
$ws.Range("K62").Value = ""

# Row 63
$ws.Range("A63").Value = "cls-obo-ncit-year-date"
$ws.Range("B63").Value = "CLS OBO NCIT Year Date"
$ws.Range("C63").Value = "Observation Category Codes#social-history, null#NCIT_C20189, null#NCIT_C41009, null#NCIT_C21514, null#NCIT_C25164, null#NCIT_C159612"
$ws.Range("D63").Value = ""
$ws.Range("E63").Value = ""
$ws.Range("F63").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G63").Value = "dateTime, Period, Timing, instant"
$ws.Range("H63").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I63").Value = "optional"
$ws.Range("J63").Value = ""
$ws.Range("K63").Value = ""

# Row 64
$ws.Range("A64").Value = "cls-obo-snomed-blood-creatinine-level"
$ws.Range("B64").Value = "CLS OBO SNOMED Blood Creatinine Level"
$ws.Range("C64").Value = "Observation Category Codes#laboratory, null#CMO_0000000, null#CMO_0000035, null#CMO_0000023, null#CMO_0000767, null#CMO_0000538"
$ws.Range("D64").Value = ""
$ws.Range("E64").Value = ""
$ws.Range("F64").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G64").Value = "dateTime, Period, Timing, instant"
$ws.Range("H64").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I64").Value = "optional"
$ws.Range("J64").Value = ""
$ws.Range("K64").Value = ""

# Row 65
$ws.Range("A65").Value = "cls-snomed-demographic-history-detail"
$ws.Range("B65").Value = "CLS SNOMED Demographic History Detail"
$ws.Range("C65").Value = "Observation Category Codes#social-history, SNOMED CT#138875005, SNOMED CT#363787002, SNOMED CT#160476009, SNOMED CT#302147001"
$ws.Range("D65").Value = ""
$ws.Range("E65").Value = ""
$ws.Range("F65").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G65").Value = "dateTime, Period, Timing, instant"
$ws.Range("H65").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I65").Value = "optional"
$ws.Range("J65").Value = ""
$ws.Range("K65").Value = ""

# Row 66
$ws.Range("A66").Value = "cls-snomed-finding-of-tobacco-use-and-exposure"
$ws.Range("B66").Value = "CLS SNOMED Finding Of Tobacco Use And Exposure"
$ws.Range("C66").Value = "Observation Category Codes#social-history, SNOMED CT#138875005, SNOMED CT#404684003, SNOMED CT#384821006, SNOMED CT#365949003, SNOMED CT#365980008"
$ws.Range("D66").Value = ""
$ws.Range("E66").Value = ""
$ws.Range("F66").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G66").Value = "dateTime, Period, Timing, instant"
$ws.Range("H66").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I66").Value = "optional"
$ws.Range("J66").Value = ""
$ws.Range("K66").Value = ""

# Row 67
$ws.Range("A67").Value = "cls-snomed-patient-sex"
$ws.Range("B67").Value = "CLS SNOMED Patient Sex"
$ws.Range("C67").Value = "Observation Category Codes#social-history, SNOMED CT#138875005, SNOMED CT#363787002, SNOMED CT#160476009, SNOMED CT#302147001, SNOMED CT#184100006"
$ws.Range("D67").Value = ""
$ws.Range("E67").Value = ""
$ws.Range("F67").Value = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G67").Value = "dateTime, Period, Timing, instant"
$ws.Range("H67").Value = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I67").Value = "optional"
$ws.Range("J67").Value = ""
$ws.Range("K67").Value = ""

# Row 68
$ws.Range("A68").Value = "myObservation"
$ws.Range("B68").Value = "My Observation Profile"
$ws.Range("C68").Value = ""
$ws.Range("D68").Value = ""
$ws.Range("E68").Value = "LOINC#85354-9"
$ws.Range("F68").Value = ""
$ws.Range("G68").Value = "dateTime, Period, Timing, instant"
$ws.Range("H68").Value = "Quantityĵ, CodeableConceptĵ"
$ws.Range("I68").Value = "optional"
$ws.Range("J68").Value = ""
$ws.Range("K68").Value = ""

# Row 69
$ws.Range("A69").Value = ""
$ws.Range("B69").Value = "My Observation Profile"
$ws.Range("C69").Value = ""
$ws.Range("D69").Value = ""
$ws.Range("E69").Value = "LOINC#8480-6"
$ws.Range("F69").Value = ""
$ws.Range("G69").Value = ""
$ws.Range("H69").Value = "Quantity"
$ws.Range("I69").Value = "optional"
$ws.Range("J69").Value = ""
$ws.Range("K69").Value = ""

# Row 70
$ws.Range("A70").Value = ""
$ws.Range("B70").Value = "My Observation Profile"
$ws.Range("C70").Value = ""
$ws.Range("D70").Value = ""
$ws.Range("E70").Value = "LOINC#8462-4"
$ws.Range("F70").Value = ""
$ws.Range("G70").Value = ""
$ws.Range("H70").Value = "Quantity"
$ws.Range("I70").Value = "optional"
$ws.Range("J70").Value = ""
$ws.Range("K70").Value = ""

# --- Remove now-unused trailing rows (old cls-l3-* / myObservation block tail) ---
$ws.Range("A71:K82").EntireRow.Delete()